$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D3").Value = "2016-02-18 03:40:44"
$wsZh.Range("G3").Value = "2016-02-18 03:41:29"

$wsDe.Range("D3").Value = "2016-02-18 03:40:57"
$wsDe.Range("G3").Value = "2016-02-18 03:41:51"
